$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 222443.22
$ws.Range("I12").Value = 164.83333
$ws.Range("J12").Value = 667000
$ws.Range("K12").Value = 164.83333
$ws.Range("L12").Value = 667000
$ws.Range("M12").Value = 5.166670000000011
$ws.Range("N12").Value = -667340
$ws.Range("H62").Value = 2704.5386
$ws.Range("I62").Value = 2114.8572
$ws.Range("K62").Value = 2114.8572
$ws.Range("M62").Value = -1490.8572
$ws.Range("H65").Value = 2704.5386
$ws.Range("I65").Value = 2114.8572
$ws.Range("K65").Value = 10574.286
$ws.Range("M65").Value = -7454.286
$ws.Range("H82").Value = 443.5
$ws.Range("I82").Value = 443.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1330.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -924.5
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 443.5
$ws.Range("I85").Value = 443.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1330.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 73.5
$ws.Range("N85").ClearContents()
$ws.Range("H97").Value = 7416.625
$ws.Range("J97").Value = 7416.625
$ws.Range("L97").Value = 22249.875
$ws.Range("N97").Value = -23241.875
$ws.Range("H100").Value = 3594.524
$ws.Range("I100").Value = 3565.4167
$ws.Range("J100").Value = 3633.3333
$ws.Range("K100").Value = 3565.4167
$ws.Range("L100").Value = 3633.3333
$ws.Range("M100").Value = -3024.4167
$ws.Range("N100").Value = -4715.3333
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H138").Value = 3245.54
$ws.Range("I138").Value = 2202.4883
$ws.Range("J138").Value = 4032.4036
$ws.Range("K138").Value = 6607.4649
$ws.Range("L138").Value = 12097.2108
$ws.Range("M138").Value = -1467.4649
$ws.Range("N138").Value = -22377.2108

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14250.142
$ws.Range("I32").Value = 12182.646
$ws.Range("J32").Value = 20969.5
$ws.Range("K32").Value = 12182.646
$ws.Range("L32").Value = 20969.5
$ws.Range("M32").Value = -11895.646
$ws.Range("N32").Value = -21543.5
$ws.Range("H61").Value = 2319.0833
$ws.Range("I61").Value = 1650.381
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 1650.381
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -1438.381
$ws.Range("N61").Value = -7424
$ws.Range("H74").Value = 1798.3572
$ws.Range("I74").Value = 1376.8096
$ws.Range("J74").Value = 3063
$ws.Range("K74").Value = 1376.8096
$ws.Range("L74").Value = 3063
$ws.Range("M74").Value = -502.8096
$ws.Range("N74").Value = -4811
$ws.Range("H77").Value = 1798.3572
$ws.Range("I77").Value = 1376.8096
$ws.Range("J77").Value = 3063
$ws.Range("K77").Value = 6884.048000000001
$ws.Range("L77").Value = 15315
$ws.Range("M77").Value = -2516.048000000001
$ws.Range("N77").Value = -24051
$ws.Range("H97").Value = 524
$ws.Range("I97").Value = 524
$ws.Range("K97").Value = 524
$ws.Range("M97").Value = -28
$ws.Range("H102").Value = 4395
$ws.Range("I102").Value = 3375.75
$ws.Range("K102").Value = 3375.75
$ws.Range("M102").Value = -1753.75
$ws.Range("H136").Value = 2319.0833
$ws.Range("I136").Value = 1650.381
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 4951.143
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -2401.143
$ws.Range("N136").Value = -26100

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 914.6
$ws.Range("I94").Value = 734.3333
$ws.Range("J94").Value = 1185
$ws.Range("K94").Value = 734.3333
$ws.Range("L94").Value = 1185
$ws.Range("M94").Value = -283.3333
$ws.Range("N94").Value = -2087
$ws.Range("H105").Value = 1657.2307
$ws.Range("I105").Value = 1451
$ws.Range("K105").Value = 1451
$ws.Range("M105").Value = 296

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2850.705
$ws.Range("I31").Value = 1862.0682
$ws.Range("J31").Value = 4130.1177
$ws.Range("K31").Value = 1862.0682
$ws.Range("L31").Value = 4130.1177
$ws.Range("M31").Value = -1567.0682
$ws.Range("N31").Value = -4720.1177
$ws.Range("H34").Value = 2850.705
$ws.Range("I34").Value = 1862.0682
$ws.Range("J34").Value = 4130.1177
$ws.Range("K34").Value = 1862.0682
$ws.Range("L34").Value = 4130.1177
$ws.Range("M34").Value = -1660.0682
$ws.Range("N34").Value = -4534.1177
$ws.Range("H58").Value = 10003469
$ws.Range("I58").Value = 1752.2307
$ws.Range("J58").Value = 20838662
$ws.Range("K58").Value = 1752.2307
$ws.Range("L58").Value = 20838662
$ws.Range("M58").Value = -1549.2307
$ws.Range("N58").Value = -20839068
$ws.Range("H97").Value = 32700
$ws.Range("J97").Value = 32700
$ws.Range("L97").Value = 32700
$ws.Range("N97").Value = -34682
$ws.Range("H132").Value = 1901.6571
$ws.Range("I132").Value = 1690.9615
$ws.Range("J132").Value = 2510.3333
$ws.Range("K132").Value = 5072.8845
$ws.Range("L132").Value = 7530.999899999999
$ws.Range("M132").Value = -2542.8845
$ws.Range("N132").Value = -12590.9999
$ws.Range("H133").Value = 25800
$ws.Range("J133").Value = 25800
$ws.Range("L133").Value = 25800
$ws.Range("N133").Value = -30860
$ws.Range("H134").Value = 2548.7173
$ws.Range("I134").Value = 1659
$ws.Range("J134").Value = 4807.231
$ws.Range("K134").Value = 4977
$ws.Range("L134").Value = 14421.693
$ws.Range("M134").Value = -2442
$ws.Range("N134").Value = -19491.693
$ws.Range("H135").Value = 18930.225
$ws.Range("J135").Value = 18930.225
$ws.Range("L135").Value = 18930.225
$ws.Range("N135").Value = -29070.225
$ws.Range("H136").Value = 10003469
$ws.Range("I136").Value = 1752.2307
$ws.Range("J136").Value = 20838662
$ws.Range("K136").Value = 5256.6921
$ws.Range("L136").Value = 62515986
$ws.Range("M136").Value = -2706.6921
$ws.Range("N136").Value = -62521086
$ws.Range("H137").Value = 76593.336
$ws.Range("J137").Value = 76593.336
$ws.Range("L137").Value = 76593.336
$ws.Range("N137").Value = -86793.336

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2714.6667
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 2857.6
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 8572.799999999999
$ws.Range("M125").Value = -1080
$ws.Range("N125").Value = -18412.8
$ws.Range("H131").Value = 1526.3265
$ws.Range("I131").Value = 3981.6667
$ws.Range("J131").Value = 1183.721
$ws.Range("K131").Value = 11945.0001
$ws.Range("L131").Value = 3551.163
$ws.Range("M131").Value = -6905.000100000001
$ws.Range("N131").Value = -13631.163

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 51122.11
$ws.Range("I52").Value = 15000
$ws.Range("J52").Value = 80019.8
$ws.Range("K52").Value = 15000
$ws.Range("L52").Value = 80019.8
$ws.Range("N52").Value = -80537.8
$ws.Range("M52").Value = -14741
$ws.Range("H97").Value = 1840.7693
$ws.Range("I97").Value = 1268.7778
$ws.Range("J97").Value = 3127.75
$ws.Range("K97").Value = 1268.7778
$ws.Range("L97").Value = 3127.75
$ws.Range("M97").Value = -772.7778000000001
$ws.Range("N97").Value = -4119.75
$ws.Range("H107").Value = 1043.1052
$ws.Range("I107").Value = 454
$ws.Range("J107").Value = 1471.5454
$ws.Range("K107").Value = 454
$ws.Range("L107").Value = 1471.5454
$ws.Range("M107").Value = 1466
$ws.Range("N107").Value = -5311.5454
$ws.Range("H140").Value = 17228.121
$ws.Range("J140").Value = 17228.121
$ws.Range("L140").Value = 17228.121
$ws.Range("N140").Value = -27588.121

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 475
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 55557856
$ws.Range("I61").Value = 66667240
$ws.Range("K61").Value = 66667240
$ws.Range("M61").Value = -66667038
$ws.Range("H113").Value = 55557856
$ws.Range("I113").Value = 66667240
$ws.Range("K113").Value = 66667240
$ws.Range("M113").Value = -66665070

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1781.2
$ws.Range("I96").Value = 1618.3334
$ws.Range("J96").Value = 1889.7778
$ws.Range("K96").Value = 1618.3334
$ws.Range("L96").Value = 1889.7778
$ws.Range("M96").Value = -245.3334
$ws.Range("N96").Value = -4635.7778
